$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.988.81"
$ws.Range("E2").Value = "  +3.24%  "
$ws.Range("D3").Value = "3.265.94"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.44"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "3.266.00"
$ws.Range("E9").Value = "  +2.84%  "
$ws.Range("E10").Value = "  +7.41%  "
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("E12").Value = "  +6.04%  "
$ws.Range("D13").Value = "3.837.70"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.57"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("D16").Value = "68.016.62"
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000169"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.35%  "
$ws.Range("D18").Value = "3.270.07"
$ws.Range("E18").Value = "  +2.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.85"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.52"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.30"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.65"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.06%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("E26").Value = "  +4.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.62"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.70"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.86"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.76%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.27"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.53%  "
$ws.Range("E35").Value = "  +4.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.21"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.851"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("E39").Value = "  +2.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.83"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +11.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.93"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("E42").Value = "  +10.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.55%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "351.15"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.35%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.687.14"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.64"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.88"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0681"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.27%  "
$ws.Range("E51").Value = "  +0.53%  "

Write-Host "Applied 91 cell updates"